$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

$newValue = "CLUTEST_04_250808"

$row = 8
while ($row -le 809) {
    $ws.Cells.Item($row, 5).Value = $newValue
    $ws.Cells.Item($row + 1, 5).Value = $newValue
    $row = $row + 10
}
